# "Updated slides with delayed mode"
#
# 1. Add the pair of default centre slide guides (horizontal @ 270pt,
#    vertical @ 360pt) that PowerPoint records in presentation.xml's
#    p:extLst/p15:sldGuideLst when guides are first shown/dragged out.
# 2. Re-cache the "last modified" date field shown on the slide master,
#    every slide layout and the notes master (26.09.2012 -> 16.10.2013).
# 3. On the two "Fulfilled / Failed" legend textboxes, change the wording
#    to "Fulfilled / Delayed" (typed as two runs) and let the auto-fit
#    textbox grow to the new (wider) cached size.

$p = $ppt.ActivePresentation
$newDate = "16.10.2013"
$ppPlaceholderDate = 16

# --- 1. slide guides ---------------------------------------------------
try {
    $guides = $p.Guides
    [void]$guides.Add(1, 270)
    [void]$guides.Add(2, 360)
} catch {
    # Guides collection unavailable in this host - nothing more we can do
    # through the object model.
}

# --- 2. date placeholders on the slide master --------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $shape.TextFrame.TextRange.Text = $newDate
    }
}

# --- date placeholders on every slide layout ----------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- date placeholder on the notes master --------------------------------
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shape = $notesMaster.Shapes.Item($i)
    if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $shape.TextFrame.TextRange.Text = $newDate
    }
}

# --- 3. "Fulfilled / Failed" -> "Fulfilled / Delayed" --------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Fulfilled / Failed") {
                $tr = $shape.TextFrame.TextRange
                $tr.Text = "Fulfilled / "
                [void]$tr.InsertAfter("Delayed")
                $shape.Width = 118.491025
            }
        }
    }
}
